# Update performance dashboard 2025-12-25 00:09 - Simplified Design v3.0
#
# Refreshes the "Pattern3-Data+News" / "deepseek-v3" row with the latest
# trading metrics in both the Summary roll-up sheet (row 12) and the
# dedicated "Pattern3-Data+News" sheet (row 2).

function Update-PerformanceRow([object]$ws, [int]$row) {

    # Currency-formatted text - these are left alone by the engine's
    # automatic type inference because of the leading currency glyph.
    $ws.Cells.Item($row, 4).Value = "¥1,005,052.00"   # D: Latest Equity
    $ws.Cells.Item($row, 5).Value = "¥+5,052.00"      # E: Total Profit

    # Percent-looking text values need a leading quote so the engine
    # keeps them as literal text instead of re-interpreting them as
    # numeric percentages.
    $ws.Cells.Item($row, 6).Value  = "'+0.51%"        # F: Total Return (%)
    $ws.Cells.Item($row, 7).Value  = "'+23.36%"       # G: Annual Return (%)

    # Genuine numeric value.
    $ws.Cells.Item($row, 8).Value  = 17.501           # H: Sharpe Ratio

    $ws.Cells.Item($row, 10).Value = "'60.0%"         # J: Win Rate (%)
    $ws.Cells.Item($row, 11).Value = "'0.1009%"       # K: Avg Daily Return (%)
    $ws.Cells.Item($row, 12).Value = "'0.0840%"       # L: Return Volatility (%)

    # Genuine numeric values.
    $ws.Cells.Item($row, 13).Value = 6                # M: Trading Days
    $ws.Cells.Item($row, 14).Value = 6                # N: Files Count

    # Digit string that must stay literal text (quote-prefixed).
    $ws.Cells.Item($row, 16).Value = "'20251224"      # P: Data Date
}

$wb = $excel.ActiveWorkbook

# Summary sheet - row 12 (Pattern3-Data+News / deepseek-v3)
$wsSummary = $wb.Worksheets.Item("Summary")
Update-PerformanceRow $wsSummary 12

# Pattern3-Data+News sheet - row 2 (same model/pattern entry)
$wsPattern3 = $wb.Worksheets.Item("Pattern3-Data+News")
Update-PerformanceRow $wsPattern3 2
